$p = $ppt.ActivePresentation

# --- Slide 9: "CustomShape 2" - first bullet paragraph ---
# Original:  "E-Barber Shop system is responsible for providing best services, ..."
# Target:    "Online-Barber Shop system is responsible for providing best services, ..."
# Achieved by typing "Online" over the leading "E" character, then re-touching the
# "-Barber " span so it becomes its own run (mirrors how PowerPoint splits runs when
# text is edited in place).
$slide9 = $p.Slides.Item(9)
$shape9 = $slide9.Shapes.Item(2)

$para1 = $shape9.TextFrame.TextRange.Paragraphs(1)
$leadChar = $para1.Characters(1, 1)
$leadChar.Text = "Online"

$para1b = $shape9.TextFrame.TextRange.Paragraphs(1)
$dashBarber = $para1b.Characters(7, 8)
$dashBarber.Text = "-Barber "


